$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "62.626.07"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = "  -1.68%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.670.44"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "  -2.17%  "

$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "552.43"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "  -2.33%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "156.37"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = "  -2.33%  "

$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = "  +0.13%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.586"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = "  -1.90%  "

$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = "  -4.31%  "

$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = "  -3.83%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.50"
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = "  -2.67%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.363"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = "  -4.44%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "3.143.66"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = "  -2.08%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "25.90"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = "  -3.74%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "62.554.80"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = "  -1.49%  "

$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = "  -3.34%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.674.97"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = "  -2.02%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "11.76"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "  -6.70%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "4.55"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "  -4.13%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "342.55"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "  -3.49%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.16"
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "  -6.51%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.998"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = "  -0.11%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "0.505"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = "  -3.44%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "63.05"
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = "  -2.42%  "

$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = "  -0.21%  "

$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = "  +0.14%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "8.07"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "  -3.81%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "1.38"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = "  +3.39%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "0.0₃0839"
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = "  -7.64%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "7.20"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = "  +0.05%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.91"
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = "  -3.40%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "162.09"
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = "  -2.49%  "

$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = "  +0.05%  "

$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = "  -2.52%  "

$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = "  -2.71%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "19.34"
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = "  -3.76%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.76"
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = "  -2.46%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "334.48"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "  -3.86%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "6.10"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "  -3.24%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.914"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = "  -6.36%  "

$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = "  -3.45%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "38.21"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = "  -1.22%  "

$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = "  -5.48%  "

$ws.Cells.Item(44, 2).NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = "Mantle"
$ws.Cells.Item(44, 3).NumberFormat = "@"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.613"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = "  -1.72%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.999"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = "  +0.10%  "

$ws.Cells.Item(46, 2).NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).NumberFormat = "@"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "19.98"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = "  -5.07%  "

$ws.Cells.Item(47, 2).NumberFormat = "@"
$ws.Cells.Item(47, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(47, 3).NumberFormat = "@"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "11.00"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = "  -0.43%  "

$ws.Cells.Item(48, 2).NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = "Hedera"
$ws.Cells.Item(48, 3).NumberFormat = "@"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0551"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = "  -5.69%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "129.11"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "  -2.15%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0963"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = "  -3.84%  "

$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = "  -4.84%  "
